$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting the existing weekly records
# (rows 24-44) down to rows 25-45, then populate the freshly inserted
# row with the new week's reading.
$ws.Rows("24:24").Insert()

$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C24").Value = "Los Lagos"
$ws.Range("D24").Value = 44421
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 100112031
$ws.Range("G24").Value = "Poroto verde"
$ws.Range("H24").Value = "Magnum"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 40
$ws.Range("K24").Value = 42000
$ws.Range("L24").Value = 42000
$ws.Range("M24").Value = 42000
$ws.Range("N24").Value = "`$/malla 25 kilos"
$ws.Range("O24").Value = "Perú"
$ws.Range("P24").Value = 1680
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"

# Append a brand new record as the last row (row 46) for the most
# recent week.
$ws.Range("A46").Value = 4
$ws.Range("B46").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C46").Value = "Los Lagos"
$ws.Range("D46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D46").Value = 44418
$ws.Range("E46").Value = 10
$ws.Range("F46").Value = 100112031
$ws.Range("G46").Value = "Poroto verde"
$ws.Range("H46").Value = "Magnum"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 50
$ws.Range("K46").Value = 40000
$ws.Range("L46").Value = 40000
$ws.Range("M46").Value = 40000
$ws.Range("N46").Value = "`$/malla 25 kilos"
$ws.Range("O46").Value = "Perú"
$ws.Range("P46").Value = 1600
$ws.Range("Q46").Value = 25
$ws.Range("R46").Value = "Hortaliza"
